$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2's flight number (reuses the now-orphaned shared string slot)
$ws.Range("A2").Value = "VN001000"

# Fill in row 4
$ws.Range("A4").Value = "VN001001"
$ws.Range("B4").Value = "VN4567"
$ws.Range("C4").Value = Get-Date -Year 2024 -Month 12 -Day 9 -Hour 0 -Minute 30 -Second 0
$ws.Range("D4").Value = 90
$ws.Range("E4").Value = "Tan Son Nhat International Airport"
$ws.Range("F4").Value = "Noi Bai International Airport"
$ws.Range("G4").Value = 1490000
$ws.Range("H4").Value = 2000000
$ws.Range("I4").Value = 50000000

# Fill in row 3
$ws.Range("A3").Value = "VN001002"
$ws.Range("B3").Value = "VN4567"
$ws.Range("C3").Value = Get-Date -Year 2024 -Month 12 -Day 9 -Hour 11 -Minute 0 -Second 0
$ws.Range("D3").Value = 90
$ws.Range("E3").Value = "Tan Son Nhat International Airport"
$ws.Range("F3").Value = "Noi Bai International Airport"
$ws.Range("G3").Value = 1490000
$ws.Range("H3").Value = 2000000
$ws.Range("I3").Value = 50000000

# Update selection to match final state
$ws.Range("B5").Select()
